# Updated for ms sql
# Replace the old order-user email with the new one across all the
# "order" sheets that reference it, then leave the UI state (selected
# cells, active sheet, column width) the way Excel would after a user
# made this edit interactively.

$wb = $excel.ActiveWorkbook

$newUser = "orderuser1@10004.escm.local"

$sheetNames = @(
    "createOrder",
    "upsizeOrder",
    "downsizeOrder",
    "upgradeOrder",
    "downgradeOrder",
    "cancleFromStableState"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A2").Value = $newUser
    $ws.Range("A2").Select() | Out-Null
}

# createOrder: widen column A so the longer address is visible, and
# leave the A2 cell selected (previously G2 was selected).
$wsCreate = $wb.Worksheets.Item("createOrder")
$wsCreate.Columns.Item(1).ColumnWidth = 35.5
$wsCreate.Range("A2").Select() | Out-Null

# upgradeOrder / downgradeOrder: selection moves from the last column to A2.
$wsUpgrade = $wb.Worksheets.Item("upgradeOrder")
$wsUpgrade.Range("A2").Select() | Out-Null

$wsDowngrade = $wb.Worksheets.Item("downgradeOrder")
$wsDowngrade.Range("A2").Select() | Out-Null

# Re-order which sheet is scrolled-to-first and which is active: the
# workbook now opens showing "upsizeOrder" as the first visible tab,
# with "cancleFromStableState" the active/selected sheet.
$wsCancel = $wb.Worksheets.Item("cancleFromStableState")
$wsCancel.Activate()

$wb.Windows.Item(1).ScrollWorkbookTabs(-17) | Out-Null
$wb.Windows.Item(1).ScrollWorkbookTabs(5) | Out-Null
